$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$type02 = "Type 02: demand in imperative form + threat in active form"
$type03 = "Type 03: threat in active form + demand in active form"
$noExt  = "No Extortion"

# Rows changing from Type 02 -> Type 03
$ws.Range("B12").Value  = $type03
$ws.Range("B23").Value  = $type03
$ws.Range("B113").Value = $type03
$ws.Range("B366").Value = $type03

# Rows changing from Type 03 -> Type 02
$ws.Range("B750").Value = $type02
$ws.Range("B774").Value = $type02
$ws.Range("B806").Value = $type02
$ws.Range("B847").Value = $type02
$ws.Range("B880").Value = $type02

# Rows changing from No Extortion -> Type 03
$ws.Range("B756").Value = $type03
$ws.Range("B765").Value = $type03
$ws.Range("B784").Value = $type03
$ws.Range("B796").Value = $type03
$ws.Range("B816").Value = $type03
$ws.Range("B819").Value = $type03
$ws.Range("B822").Value = $type03
$ws.Range("B863").Value = $type03
$ws.Range("B895").Value = $type03
$ws.Range("B908").Value = $type03
$ws.Range("B916").Value = $type03
$ws.Range("B919").Value = $type03
$ws.Range("B925").Value = $type03
$ws.Range("B971").Value = $type03

# Row changing from Type 03 -> No Extortion
$ws.Range("B969").Value = $noExt
